# Update of league bases - swap the data (all columns except the id in
# column A) between the listed pairs of rows. Column A (the sequential
# record id) stays attached to its original row; every other field
# (match id, teams, scores, odds, etc.) moves with its match record to
# the other row in the pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsA = @(18, 36, 38, 58, 63)
$rowsB = @(19, 37, 39, 59, 64)

for ($i = 0; $i -lt $rowsA.Count; $i++) {
    $row1 = $rowsA[$i]
    $row2 = $rowsB[$i]

    $range1 = $ws.Range("B${row1}:AD${row1}")
    $range2 = $ws.Range("B${row2}:AD${row2}")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value = $values2
    $range2.Value = $values1
}
